$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '98.165.57'
$ws.Range("E2").Value = '  +4.50%  '

# Row 3
$ws.Range("D3").Value = '3.356.24'
$ws.Range("E3").Value = '  +9.64%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '623.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.87%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.386'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.46%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").Value = '3.355.39'
$ws.Range("E10").Value = '  +9.73%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.800'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '

# Row 12
$ws.Range("E12").Value = '  +2.53%  '

# Row 13
$ws.Range("D13").Value = '97.916.23'
$ws.Range("E13").Value = '  +4.58%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.66%  '

# Row 15
$ws.Range("E15").Value = '  +3.52%  '

# Row 16
$ws.Range("D16").Value = '3.959.92'
$ws.Range("E16").Value = '  +9.09%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.98%  '

# Row 18
$ws.Range("D18").Value = '3.356.20'
$ws.Range("E18").Value = '  +10.30%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.26%  '

# Row 20
$ws.Range("E20").Value = '  +5.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '486.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.75%  '

# Row 22
$ws.Range("E22").Value = '  +4.13%  '

# Row 23
$ws.Range("E23").Value = '  +11.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.78%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.65%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.11%  '

# Row 28
$ws.Range("D28").Value = '3.531.89'
$ws.Range("E28").Value = '  +9.59%  '

# Row 29
$ws.Range("E29").Value = '  +0.09%  '

# Row 30
$ws.Range("E30").Value = '  +5.47%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.246'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.61%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.124'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.28%  '

# Row 34
$ws.Range("E34").Value = '  +3.89%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.69%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.74%  '

# Row 37
$ws.Range("E37").Value = '  -1.52%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '515.64'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.69%  '

# Row 39
$ws.Range("E39").Value = '  +4.01%  '

# Row 40
$ws.Range("E40").Value = '  +3.77%  '

# Row 41
$ws.Range("E41").Value = '  +4.57%  '

# Row 42
$ws.Range("E42").Value = '  +2.53%  '

# Row 43
$ws.Range("E43").Value = '  -3.43%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.88%  '

# Row 46
$ws.Range("E46").Value = '  +18.08%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '161.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.76%  '

# Row 48
$ws.Range("E48").Value = '  +8.03%  '

# Row 49
$ws.Range("E49").Value = '  +8.97%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.89%  '
